$d = $word.ActiveDocument

# 1. "A collision is an event..." paragraph: insert "also" so that
#    "collisions can occur between repelling particles." becomes
#    "collisions can also occur between repelling particles."
#    (only the occurrence immediately followed by "Particles of equal charge
#    sign" should change - the sentence is duplicated elsewhere in the doc).
$d.Content.Find.Execute(
    "collisions can occur between repelling particles. Particles of equal charge sign",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "collisions can also occur between repelling particles. Particles of equal charge sign",
    2) | Out-Null

# 2. "the electron is liberated and the atom is" -> "the electron is liber and the atom is"
$d.Content.Find.Execute(
    "the electron is liberated and the atom is",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the electron is liber and the atom is",
    2) | Out-Null

# 3. "infrequently relative to the occurrence" -> "infrequent relative to the occurrence"
$d.Content.Find.Execute(
    "infrequently relative to the occurrence",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "infrequent relative to the occurrence",
    2) | Out-Null

# 4. Remove the duplicated "Elastic collision" bullet item directly following
#    "Inelastic collision" in the interactions list (keep the other copy).
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Elastic collision`r") {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text -eq "Inelastic collision`r") {
            $p.Range.Delete()
            break
        }
    }
}
